# Insert a new "Demo" slide between "CRUD in the ORM" (slide 11) and
# "Summary" (slide 12), so it becomes the new slide 12 and pushes
# "Summary" / "Acknowledgements / Contributions" down by one position.

$p = $ppt.ActivePresentation

# 1. Insert a new slide at position 12 using the "Title and Content" layout
#    (same layout already used by several other slides in this deck).
$s = $p.Slides.Add(12, 2)

# 2. Give the new slide the same dark "gears" background picture used on
#    the title slide (slide 1, shape 2 = "Picture 3"), then push it to the
#    back of the z-order so the title/body placeholders sit on top of it.
$bgPic = $p.Slides.Item(1).Shapes.Item(2)
$bgPic.Copy()
$pastedShapes = $s.Shapes.Paste()
$newPic = $pastedShapes.Item(1)
$newPic.ZOrder(1)

# 3. Set the title placeholder text.
$title = $s.Shapes.Item(2)
$title.TextFrame.TextRange.Text = "Demo"

# 4. Set the body placeholder text (typed as two chunks, like autocomplete
#    splitting off the "https://" prefix) and rename it to match the
#    target shape name, then add a trailing blank paragraph.
$body = $s.Shapes.Item(3)
$body.Name = "Text Placeholder 2"
$bodyTextRange = $body.TextFrame.TextRange
$bodyTextRange.Text = "https://"
$bodyTextRange.InsertAfter("www.dj4e.com/lectures/DJ-03-Model-Single.txt")
$body.TextFrame.TextRange.InsertAfter([char]13)
